$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.793.12'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('D3').Value = '2.673.93'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = "'603.10"
$ws.Range('E5').Value = '  -1.52%  '
$ws.Range('D6').Value = "'157.83"
$ws.Range('E6').Value = '  -1.35%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  +5.18%  '
$ws.Range('D9').Value = "'0.130"
$ws.Range('E9').Value = '  +2.95%  '
$ws.Range('D10').Value = "'0.404"
$ws.Range('E10').Value = '  -1.01%  '
$ws.Range('E11').Value = '  -2.81%  '
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('D13').Value = "'29.61"
$ws.Range('E13').Value = '  -2.59%  '
$ws.Range('D14').Value = "'0.0000201"
$ws.Range('E14').Value = '  -6.56%  '
$ws.Range('D15').Value = '3.154.49'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').Value = '65.562.66'
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').Value = '2.673.23'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').Value = "'12.82"
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').Value = "'4.83"
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').Value = "'7.66"
$ws.Range('E20').Value = '  +2.20%  '
$ws.Range('D21').Value = "'352.90"
$ws.Range('E21').Value = '  -2.85%  '
$ws.Range('D22').Value = "'0.999"
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').Value = "'69.73"
$ws.Range('E23').Value = '  -0.51%  '
$ws.Range('D24').Value = "'0.0000112"
$ws.Range('E24').Value = '  +2.27%  '
$ws.Range('D25').Value = "'9.86"
$ws.Range('E25').Value = '  +1.92%  '
$ws.Range('E26').Value = '  -4.07%  '
$ws.Range('D27').Value = "'0.168"
$ws.Range('E27').Value = '  -3.56%  '
$ws.Range('E28').Value = '  -4.08%  '
$ws.Range('D29').Value = "'8.12"
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = "'2.16"
$ws.Range('E31').Value = '  -3.00%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = "'532.52"
$ws.Range('E32').Value = '  -2.16%  '
$ws.Range('D33').Value = "'1.77"
$ws.Range('E33').Value = '  -2.35%  '
$ws.Range('D34').Value = "'6.59"
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').Value = "'5.55"
$ws.Range('E35').Value = '  -1.18%  '
$ws.Range('D36').Value = "'0.425"
$ws.Range('E36').Value = '  -2.72%  '
$ws.Range('D37').Value = "'20.48"
$ws.Range('E37').Value = '  -1.75%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = "'159.61"
$ws.Range('E38').Value = '  -1.95%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').Value = "'0.999"
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('D40').Value = "'1.95"
$ws.Range('E40').Value = '  -4.14%  '
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('D42').Value = "'42.69"
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('D43').Value = "'165.33"
$ws.Range('E43').Value = '  -3.11%  '
$ws.Range('D44').Value = "'4.11"
$ws.Range('E44').Value = '  -2.51%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = "'2.33"
$ws.Range('E45').Value = '  -2.85%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = "'0.0614"
$ws.Range('E46').Value = '  -0.93%  '
$ws.Range('D47').Value = "'23.15"
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('D48').Value = "'0.0261"
$ws.Range('E48').Value = '  -1.94%  '
$ws.Range('D49').Value = "'0.645"
$ws.Range('E49').Value = '  -2.97%  '
$ws.Range('D50').Value = "'0.102"
$ws.Range('E50').Value = '  +3.27%  '
$ws.Range('D51').Value = "'20.39"
$ws.Range('E51').Value = '  +1.60%  '
